$d = $word.ActiveDocument

# Change 1
$d.Content.Find.Execute("enlarged to full page size, in the supplementary materials.", $true, $false, $false, $false, $false, $true, 1, $false, "enlarged to span each page, in the Supplementary Materials.", 2) | Out-Null

# Change 2
$d.Content.Find.Execute("The top rows of the revised figures 3-7 report the mean sampling rates and in the top rows of figures S7, S9-S11 and S14 the ranks of chosen options for all the models. We plot these alongside those of participants for direct comparison of whether the models can reproduce participants’ performance. ", $true, $false, $false, $false, $false, $true, 1, $false, "The top rows of the revised Figures 3-7 report the mean sampling rates and in the top rows of Figures S7, S9-S11 and S14 the mean ranks of chosen options for all the models. We plot these measures alongside participant behaviour for direct comparison of whether the models reproduce participants’ performance. ", 2) | Out-Null

# Change 3
$d.Content.Find.Execute("It can be seen that all models reasonably reproduce participants’ sampling rates and the Cost to Sample and Biased Prior models (but not the Cut Off heuristic) closely reproduce participants’ mean rank of chosen option. ", $true, $false, $false, $false, $false, $true, 1, $false, "It can be seen that all models reasonably reproduce participants’ sampling rates. It can also be seen that, at the level of mean ranks of chosen prices, the models do diverge. Specifically, the Cut Off heuristic fails to reproduce participants’ performance while the Cost to Sample and Biased Prior models can. ", 2) | Out-Null

# Change 4
$d.Content.Find.Execute("In our revised General Discussion beginning line 697 we", $true, $false, $false, $false, $false, $true, 1, $false, "In our revised General Discussion beginning line 700 we", 2) | Out-Null

# Change 5
$d.Content.Find.Execute("Indeed, all our models well-predicted participants’ mean sampling rates (Figures 3-7), participants’ mean rank of chosen options (Figures S7, S9-S11, S14), except for the Cut Off heuristic, which obtained much lower ranks than participants. Our supplementary analyses of the large sample in Study 2 show individual differences in participant sampling rates were well-predicted by all three models (Figure S12). ", $true, $false, $false, $false, $false, $true, 1, $false, "Indeed, all our models accurately predicted participants’ mean sampling rates (Figures 3-7). The models did however diverge in their predictions of participants’ mean rank of chosen options (Figures S7, S9-S11, S14). Here, the Cut Off heuristic was unable to obtain similar levels as participants, while Biased Prior and Cost to Sample models could. Our supplementary analyses of the large sample in Study 2 show individual differences in participant sampling rates were highly correlated with sampling rates predicted by all three models (Figure S12). ", 2) | Out-Null

# Change 6
$d.Content.Find.Execute("the Cost to Sample model best fit a remarkable share", $true, $false, $false, $false, $false, $true, 1, $false, "the Cost to Sample model best fitted a remarkable share", 2) | Out-Null

# Change 7
$d.Content.Find.Execute("We write in the General Discussion:", $true, $false, $false, $false, $false, $true, 1, $false, "We write in the General Discussion starting line 799:", 2) | Out-Null

# Change 8
$d.Content.Find.Execute(" (the cut off model) ", $true, $false, $false, $false, $false, $true, 1, $false, " (the Cut Off heuristic) ", 2) | Out-Null

# Change 9
$d.Content.Find.Execute("or use a simpler heuristic to approximate the choice threshold / value of sampling again. Though we note that already", $true, $false, $false, $false, $false, $true, 1, $false, "or create some simpler heuristic to approximate the choice threshold / value of sampling again. We note that already", 2) | Out-Null

# Change 10
$d.Content.Find.Execute("real losses for participants. We have already proposed above an interesting theoretical possibility that biases like Biased Prior strategies might have an adaptive function, so long as they produce near-optimal performance. Indeed, within the narrow range of sequence lengths and domains (i.e., smartphones prices) that we have examined here,", $true, $false, $false, $false, $false, $true, 1, $false, "real losses for people confronted with real optimal stopping problems. We have already proposed above an interesting theoretical possibility that biases like Biased Prior strategies might have an adaptive function, so long as they can maintain near-optimal performance. Indeed, within the narrow range of sequence lengths and domains (i.e., smartphone prices) that we have examined here,", 2) | Out-Null

# Change 11
$d.Content.Find.Execute("ever larger undersampling biases cannot be answered directly by our data and would benefit from more direct investigations.", $true, $false, $false, $false, $false, $true, 1, $false, "ever larger biases cannot be answered directly by our data and would benefit from more targeted investigations.", 2) | Out-Null
